$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").Value = '6260261fcb5c2af8bf638441'
$ws.Range("B8").Value = 'Song Wendy'
$ws.Range("C8").Value = 'North Korea'
$ws.Range("D8").Value = 'songwendy@g.c'
$ws.Range("G8").Value = 'Single'
$ws.Range("H8").Value = 'No Reason'
$ws.Range("I8").NumberFormat = "@"
$ws.Range("I8").Value = '3434'
$ws.Range("J8").Value = 'Pending'
$ws.Range("K8").Value = 'image_1650468383148_gcash-merged.pdf'
$ws.Range("L8").Value = 'April 20th 2022'
$ws.Range("E8").Value = 435
$ws.Range("F8").Value = 98
$ws.Range("M8").Value = 0

# Row 9
$ws.Range("A9").Value = '6260b3d7979a8f72c0e5f8ee'
$ws.Range("B9").Value = 'Bae Irene'
$ws.Range("C9").Value = 'South Korea'
$ws.Range("D9").Value = 'baeirene@g.c'
$ws.Range("G9").Value = 'Single'
$ws.Range("H9").Value = 'No Reason'
$ws.Range("I9").Value = '25366523fdfdf'
$ws.Range("J9").Value = 'Pending'
$ws.Range("K9").Value = 'image_1650504663047_gcash-merged.pdf'
$ws.Range("L9").Value = 'April 21st 2022'
$ws.Range("E9").Value = 43090
$ws.Range("F9").Value = 29
$ws.Range("M9").Value = 0

# Row 10
$ws.Range("A10").Value = '6260b48dcdeded58dd9be7c4'
$ws.Range("B10").Value = 'Juliet Mediona Nicanor'
$ws.Range("C10").Value = 'Rosario Village, Botong Francisco Ave'
$ws.Range("D10").Value = 'julietnicanor1996@gmail.com'
$ws.Range("G10").Value = 'Single'
$ws.Range("H10").Value = 'No Reason'
$ws.Range("I10").Value = 'dsds'
$ws.Range("J10").Value = 'Pending'
$ws.Range("K10").Value = 'image_1650504845890_gcash-merged.pdf'
$ws.Range("L10").Value = 'April 21st 2022'
$ws.Range("E10").Value = 639395029337
$ws.Range("F10").Value = 23
$ws.Range("M10").Value = 0

# Row 11
$ws.Range("A11").Value = '6260b9387982862243e193d9'
$ws.Range("B11").Value = 'Juliet Mediona Nicanor'
$ws.Range("C11").Value = 'Rosario Village, Botong Francisco Ave'
$ws.Range("D11").Value = 'julietnicanor1996@gmail.com'
$ws.Range("G11").Value = 'Single'
$ws.Range("H11").Value = 'dsds'
$ws.Range("I11").Value = 'Hellow'
$ws.Range("J11").Value = 'Pending'
$ws.Range("K11").Value = 'image_1650506040165_gcash-merged.pdf'
$ws.Range("L11").Value = 'April 21st 2022'
$ws.Range("E11").Value = 639395029337
$ws.Range("F11").Value = 23
$ws.Range("M11").Value = 0

# Row 12
$ws.Range("A12").Value = '6260ba75443489731d5dd848'
$ws.Range("B12").Value = 'Juliet Mediona Nicanor'
$ws.Range("C12").Value = 'Rosario Village, Botong Francisco Ave'
$ws.Range("D12").Value = 'julietnicanor1996@gmail.com'
$ws.Range("G12").Value = 'dsd'
$ws.Range("H12").Value = 'No Reason'
$ws.Range("I12").Value = 'sas'
$ws.Range("J12").Value = 'Pending'
$ws.Range("K12").Value = 'image_1650506357140_gcash-merged.pdf'
$ws.Range("L12").Value = 'April 21st 2022'
$ws.Range("E12").Value = 639395029337
$ws.Range("F12").Value = 23
$ws.Range("M12").Value = 0

# Row 13
$ws.Range("A13").Value = '62613a069b232834da3435b1'
$ws.Range("B13").Value = 'Juliet Mediona Nicanor'
$ws.Range("C13").Value = 'Rosario Village, Botong Francisco Ave'
$ws.Range("D13").Value = 'julietnicanor1996@gmail.com'
$ws.Range("G13").Value = 'Single'
$ws.Range("H13").Value = 'dsds'
$ws.Range("I13").Value = 'sasa'
$ws.Range("J13").Value = 'Pending'
$ws.Range("K13").Value = 'image_1650539014655_gcash-merged.pdf'
$ws.Range("L13").Value = 'April 21st 2022'
$ws.Range("E13").Value = 639395029337
$ws.Range("F13").Value = 23
$ws.Range("M13").Value = 0

# Row 14
$ws.Range("A14").Value = '62613a6fc4a527c4d9811a56'
$ws.Range("B14").Value = 'Juliet Mediona Nicanor'
$ws.Range("C14").Value = 'Rosario Village, Botong Francisco Ave'
$ws.Range("D14").Value = 'julietnicanor1996@gmail.com'
$ws.Range("G14").Value = 'Single'
$ws.Range("H14").Value = 'No Reason'
$ws.Range("I14").Value = 'dsds'
$ws.Range("J14").Value = 'Pending'
$ws.Range("K14").Value = 'image_1650539119186_gcash-merged.pdf'
$ws.Range("L14").Value = 'April 21st 2022'
$ws.Range("E14").Value = 639395029337
$ws.Range("F14").Value = 23
$ws.Range("M14").Value = 0

# Row 15
$ws.Range("A15").Value = '62613bc3d1a757dc31a3791c'
$ws.Range("B15").Value = 'Juliet Mediona Nicanor'
$ws.Range("C15").Value = 'Rosario Village, Botong Francisco Ave'
$ws.Range("D15").Value = 'julietnicanor1996@gmail.com'
$ws.Range("G15").Value = 'Single'
$ws.Range("H15").Value = 'No Reason'
$ws.Range("I15").Value = 'fdfd'
$ws.Range("J15").Value = 'Pending'
$ws.Range("K15").Value = 'image_1650539459882_gcash-merged.pdf'
$ws.Range("L15").Value = 'April 21st 2022'
$ws.Range("E15").Value = 639395029337
$ws.Range("F15").Value = 23
$ws.Range("M15").Value = 0

# Row 16
$ws.Range("A16").Value = '62613cad6ad815c490fa8b18'
$ws.Range("B16").Value = 'Juliet Mediona Nicanor'
$ws.Range("C16").Value = 'Rosario Village, Botong Francisco Ave'
$ws.Range("D16").Value = 'julietnicanor1996@gmail.com'
$ws.Range("G16").Value = 'dsd'
$ws.Range("H16").Value = 'No Reason'
$ws.Range("I16").Value = 'Hellow'
$ws.Range("J16").Value = 'Pending'
$ws.Range("K16").Value = 'image_1650539693584_gcash-merged.pdf'
$ws.Range("L16").Value = 'April 21st 2022'
$ws.Range("E16").Value = 639395029337
$ws.Range("F16").Value = 23
$ws.Range("M16").Value = 0

# Row 17
$ws.Range("A17").Value = '626227a8f9a6dab50754a5a5'
$ws.Range("B17").Value = '  Duterte, Rodrigo    '
$ws.Range("C17").Value = 'Davao'
$ws.Range("D17").Value = 'duterts@g.c'
$ws.Range("G17").Value = 'Married'
$ws.Range("H17").Value = 'No Reason'
$ws.Range("I17").NumberFormat = "@"
$ws.Range("I17").Value = '94749324'
$ws.Range("J17").Value = 'Pending'
$ws.Range("K17").Value = 'image_1650599848104_gcash-merged.pdf'
$ws.Range("L17").Value = 'April 22nd 2022'
$ws.Range("E17").Value = 49249
$ws.Range("F17").Value = 79
$ws.Range("M17").Value = 0

# Row 18
$ws.Range("A18").Value = '626233aac2b270ef48934036'
$ws.Range("B18").Value = ' Manoban, Chittip  '
$ws.Range("C18").Value = 'Thailand'
$ws.Range("D18").Value = 'chits@g.c'
$ws.Range("G18").Value = 'Married'
$ws.Range("H18").Value = 'No Reason'
$ws.Range("I18").Value = 'r454854d'
$ws.Range("J18").Value = 'Pending'
$ws.Range("K18").Value = 'image_1650602922301_gcash-merged.pdf'
$ws.Range("L18").Value = 'April 22nd 2022'
$ws.Range("E18").Value = 48926489349837
$ws.Range("F18").Value = 78
$ws.Range("M18").Value = 0

# Row 19
$ws.Range("A19").Value = '626238453e5f50cf8f9a50cd'
$ws.Range("B19").Value = ' thelma m. nicanor  '
$ws.Range("C19").Value = 'Rosario Village'
$ws.Range("D19").Value = 'thelmanicanor@gmail.com'
$ws.Range("G19").Value = 'Married'
$ws.Range("H19").Value = 'Identification'
$ws.Range("I19").NumberFormat = "@"
$ws.Range("I19").Value = '567678678754545'
$ws.Range("J19").Value = 'Pending'
$ws.Range("K19").Value = 'image_1650604101967_gcash-merged.pdf'
$ws.Range("L19").Value = 'April 22nd 2022'
$ws.Range("E19").Value = 65656768769789
$ws.Range("F19").Value = 59
$ws.Range("M19").Value = 0
